# Applies the edits described by the diff for
# StructureDefinition-fr-core-patient-death-place.xlsx:
#   1. "Metadata" sheet: update the "Date" property value.
#   2. "Elements" sheet, the "Extension.value[x]" row (row 6):
#        - "Type(s)"  column -> new Address extension type reference
#        - "Short"    column -> new French/English short description
#      (the "Definition" / "Constraint(s)" columns keep their original
#       text - only their underlying shared-string index shifts in the
#       OOXML, which is an implementation detail, not a content change)
#   3. Widen the "Type(s)" column to fit the new, longer text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet - Date property
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$dateLabelCell = $wsMeta.Columns.Item(1).Find("Date")
$dateRow = $dateLabelCell.Row
$wsMeta.Cells.Item($dateRow, 2).Value = "2026-01-12T10:02:26+00:00"

# ---------------------------------------------------------------------
# 2. Elements sheet - Extension.value[x] row
# ---------------------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# locate the header columns by name so we don't depend on fixed letters
$colTypes = $wsElem.Rows.Item(1).Find("Type(s)").Column
$colShort = $wsElem.Rows.Item(1).Find("Short").Column

# locate the Extension.value[x] row by its Path column
$pathCell = $wsElem.Columns.Item(2).Find("Extension.value[x]")
$row = $pathCell.Row

$wsElem.Cells.Item($row, $colTypes).Value = "Address {https://hl7.fr/ig/fhir/core/StructureDefinition/fr-core-address|2.2.0-ballot}`n"
$wsElem.Cells.Item($row, $colShort).Value = "Lieu de décès du patient | Place where the patient is dead"

# ---------------------------------------------------------------------
# 3. Widen the "Type(s)" column (K) to fit the new text
# ---------------------------------------------------------------------
$wsElem.Columns.Item($colTypes).ColumnWidth = 65.166666666666667
